$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the jenjang value in G2 (was "Ahli Pertama" -> "Pertama")
$ws.Range("G2").Value = "Pertama"

# G3 stays "Mahir" (unchanged text, but shared-string slot changes under the hood)
$ws.Range("G3").Value = "Mahir"

# Update the selected cell in the sheet view
$ws.Range("M12").Select()
